$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove old row 2 content and column E (width 108) ---
$ws.Columns.Item(5).Delete()

# --- Row 2: new work-experience entry (RedPapaz) ---
$ws.Range("A2").Value = "Asesora Científica externa"
$ws.Range("B2").Value = "Sep. - Oct. 2022"
$ws.Range("C2").Value = "\href{https://www.redpapaz.org/}{ONG RedPapaz}"
$ws.Range("D2").Value = "Bogotá, Colombia"

# --- Row 3: new work-experience entry (Protect Children) ---
$ws.Range("A3").Value = "Consultora científica internacional"
$ws.Range("B3").Value = "Jul. - Ago. 2021"
$ws.Range("C3").Value = "\href{https://www.suojellaanlapsia.fi/en}{Protect Children}"
$ws.Range("D3").Value = "Helsink, Finlandia"

# --- Row 4: blank formatted cell ---
$ws.Range("C4").Value = ""

# --- Apply formatting: left/top align + wrap text ---
$fmtRange = $ws.Range("A2:D3")
$fmtRange.HorizontalAlignment = -4131   # xlLeft
$fmtRange.VerticalAlignment = -4160     # xlTop
$fmtRange.WrapText = $true

$ws.Range("C4").HorizontalAlignment = -4131
$ws.Range("C4").VerticalAlignment = -4160
$ws.Range("C4").WrapText = $true

# --- Row heights ---
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 30

# --- Selection ---
$ws.Range("E1").Select()
